# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
# Both sheets carry the same source data, and the same seven rows changed
# in each of them.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 1046
    "F6"  = 2848
    "F7"  = 40
    "F8"  = 1858
    "F10" = 91
    "F11" = 671
    "F13" = 25
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
